$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing "root folder" note to explicitly mention the root path is given below.
$ws.Range("A2").Value = "The path must be the remaining path after the root folder, which is:"

# Insert a new row above the "Pol0_45_90_135" row (currently row 4) to hold the explicit root path,
# pushing the sample-subfolder row down to row 5.
$ws.Rows("4").Insert()

# Fill the new row 3 with the explicit root path and merge it across A3:F3, matching the style of row 2.
$ws.Range("A3").Value = "/home/masoud/Documents/four-polar/fourPolar-io/target/test-classes/fr/fresnel/fourPolar/io/imageSet/acquisition/sample/finders/excel"
$ws.Range("A3:F3").Merge()
